# Add a new "Pallets" sheet, update the "Cajas" sheet (remove the LPNs
# that are now grouped into the pallet, keep/update the rest) and update
# the "ASN" summary sheet to include one row per "batch" written (including
# the new pallet row). The "detalle" sheet is left untouched.

$wb = $excel.ActiveWorkbook

# --- 1. Create the new "Pallets" sheet as the first sheet in the workbook ---
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "Pallets"

# NOTE: worksheet variables captured before a structural change (such as
# adding a sheet, which shifts every other sheet's position) can end up
# pointing at the wrong sheet afterwards. So every worksheet reference we
# need is (re-)fetched by name *after* the sheet has been added.
$pallets = $wb.Worksheets.Item("Pallets")
$cajas = $wb.Worksheets.Item("Cajas")
$asn = $wb.Worksheets.Item("ASN")

# --- Header row for "Pallets" ---
$palletsHeaders = @("Pallet", "LPN", "Peso (kg)", "Alto (cm)", "Largo (cm)", "Ancho (cm)")
for ($col = 1; $col -le $palletsHeaders.Length; $col++) {
    $cell = $pallets.Cells.Item(1, $col)
    $cell.Value = $palletsHeaders[$col - 1]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

$palletsRows = @(
    @("Pallet1", "SAL0000004478", 2, 5, 4, 23),
    @("Pallet1", "SAL0000004492", 2, 5, 4, 23),
    @("Pallet1", "SAL0000004494", 2, 5, 4, 23)
)
for ($r = 0; $r -lt $palletsRows.Length; $r++) {
    $rowData = $palletsRows[$r]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $pallets.Cells.Item($r + 2, $c + 1).Value = $rowData[$c]
    }
}

# --- 2. Update the "Cajas" sheet: only the boxes not placed on the pallet
#        remain; clear the old rows first, then write the new values ---
$cajas.Range("A2:F7").ClearContents()

$cajasRows = @(
    @("SAL0000004491", 2, "GSP 3", 20, 56, 40),
    @("SAL0000004493", 22, "GSP 2", 31, 40, 31),
    @("SAL0000004528", 3, "GSP 3", 20, 56, 40)
)
for ($r = 0; $r -lt $cajasRows.Length; $r++) {
    $rowData = $cajasRows[$r]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $cajas.Cells.Item($r + 2, $c + 1).Value = $rowData[$c]
    }
}

# --- 3. Update the "ASN" sheet: recompute the summary, now with one row
#        per shipped unit (Unidades = 1) plus the new "Pallet" type row ---
$asn.Range("A2:F2").ClearContents()

$asnRows = @(
    @("GSP 3", 1, 2, 20, 40, 56),
    @("GSP 2", 1, 22, 31, 31, 40),
    @("GSP 3", 1, 3, 20, 40, 56),
    @("Pallet", 1, 2, 5, 23, 4)
)
for ($r = 0; $r -lt $asnRows.Length; $r++) {
    $rowData = $asnRows[$r]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $asn.Cells.Item($r + 2, $c + 1).Value = $rowData[$c]
    }
}

# "detalle" sheet is left as-is.
